$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 551
$ws.Range("I2").Value = 1481
$ws.Range("J2").Value = 6162
$ws.Range("K2").Value = 28
$ws.Range("L2").Value = 1673
$ws.Range("N2").Value = 1086
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 23
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 90
$ws.Range("S2").Value = 654
$ws.Range("T2").Value = 1146
$ws.Range("U2").Value = 77
$ws.Range("V2").Value = 9488
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 9695
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 141
$ws.Range("AA2").Value = 63
